$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "01-07-2021"
$ws.Range("B46").Value = 2482
$ws.Range("C46").Value = 856
$ws.Range("D46").Value = 1390
$ws.Range("E46").Value = 87
$ws.Range("F46").Value = 148
